$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Formula = "=_xlfn.STDEV.S(C2:C10)"
$ws.Range("N3").Formula = "=_xlfn.STDEV.S(D2:D10)"
$ws.Range("N4").Formula = "=_xlfn.STDEV.S(E2:E10)"
$ws.Range("N5").Formula = "=_xlfn.STDEV.S(F2:F10)"

$ws.Range("N6").Select()
